$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (Yats / Yats Fishers) that are no longer present
$ws.Rows.Item(55).Delete()
$ws.Rows.Item(55).Delete()

# Rewrite the data rows (2-54) with the refreshed restaurant listing
$ws.Cells.Item(2, 1).Value = 22
$ws.Cells.Item(2, 3).Value = 'A2Z Cafe (Inside and patio dining or Carry-out to Curbside)'
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = 4.7
$ws.Cells.Item(2, 6).Value = 479
$ws.Cells.Item(3, 1).Value = 12
$ws.Cells.Item(3, 3).Value = 'Aristocrat Pub & Restaurant'
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(3, 5).Value = 4.5
$ws.Cells.Item(3, 6).Value = 1182
$ws.Cells.Item(4, 1).Value = 56
$ws.Cells.Item(4, 3).Value = 'Arni''s Restaurant - Indianapolis'
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(4, 5).Value = 4.4
$ws.Cells.Item(4, 6).Value = 794
$ws.Cells.Item(5, 1).Value = 23
$ws.Cells.Item(5, 3).Value = 'Axum Ethiopian Restaurant'
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 4.7
$ws.Cells.Item(5, 6).Value = 413
$ws.Cells.Item(6, 1).Value = 48
$ws.Cells.Item(6, 3).Value = 'BRU Burger Bar'
$ws.Cells.Item(6, 4).Value = 2
$ws.Cells.Item(6, 5).Value = 4.6
$ws.Cells.Item(6, 6).Value = 4150
$ws.Cells.Item(7, 1).Value = 10
$ws.Cells.Item(7, 3).Value = 'Bluebeard'
$ws.Cells.Item(7, 4).Value = 3
$ws.Cells.Item(7, 5).Value = 4.7
$ws.Cells.Item(7, 6).Value = 1366
$ws.Cells.Item(8, 1).Value = 18
$ws.Cells.Item(8, 3).Value = 'Bonefish Grill'
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(8, 5).Value = 4.5
$ws.Cells.Item(8, 6).Value = 1081
$ws.Cells.Item(9, 1).Value = 30
$ws.Cells.Item(9, 3).Value = 'Bosphorus Istanbul Cafe'
$ws.Cells.Item(9, 4).Value = 2
$ws.Cells.Item(9, 5).Value = 4.5
$ws.Cells.Item(9, 6).Value = 1271
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 3).Value = 'Burritos & Beer Restaurant, LLC'
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 4.7
$ws.Cells.Item(10, 6).Value = 329
$ws.Cells.Item(11, 1).Value = 25
$ws.Cells.Item(11, 3).Value = 'Canal Bistro'
$ws.Cells.Item(11, 4).Value = 2
$ws.Cells.Item(11, 5).Value = 4.6
$ws.Cells.Item(11, 6).Value = 890
$ws.Cells.Item(12, 1).Value = 0
$ws.Cells.Item(12, 3).Value = 'Charleston''s Restaurant'
$ws.Cells.Item(12, 4).Value = 2
$ws.Cells.Item(12, 5).Value = 4.5
$ws.Cells.Item(12, 6).Value = 1457
$ws.Cells.Item(13, 1).Value = 55
$ws.Cells.Item(13, 3).Value = 'City Barbeque and Catering'
$ws.Cells.Item(13, 4).Value = 2
$ws.Cells.Item(13, 5).Value = 4.5
$ws.Cells.Item(13, 6).Value = 1413
$ws.Cells.Item(14, 1).Value = 4
$ws.Cells.Item(14, 3).Value = 'Cooper''s Hawk Winery & Restaurant'
$ws.Cells.Item(14, 4).Value = 2
$ws.Cells.Item(14, 5).Value = 4.6
$ws.Cells.Item(14, 6).Value = 1498
$ws.Cells.Item(15, 1).Value = 44
$ws.Cells.Item(15, 3).Value = 'Courses Restaurant'
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(15, 5).Value = 4.6
$ws.Cells.Item(15, 6).Value = 38
$ws.Cells.Item(16, 1).Value = 29
$ws.Cells.Item(16, 3).Value = 'Cracker Barrel Old Country Store'
$ws.Cells.Item(16, 4).Value = 2
$ws.Cells.Item(16, 5).Value = 4.4
$ws.Cells.Item(16, 6).Value = 2750
$ws.Cells.Item(17, 1).Value = 35
$ws.Cells.Item(17, 3).Value = 'Fire by the Monon'
$ws.Cells.Item(17, 4).Value = 2
$ws.Cells.Item(17, 5).Value = 4.6
$ws.Cells.Item(17, 6).Value = 906
$ws.Cells.Item(18, 1).Value = 37
$ws.Cells.Item(18, 3).Value = 'First Watch'
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = 4.6
$ws.Cells.Item(18, 6).Value = 396
$ws.Cells.Item(19, 1).Value = 47
$ws.Cells.Item(19, 3).Value = 'Flatwater'
$ws.Cells.Item(19, 4).Value = 2
$ws.Cells.Item(19, 5).Value = 4.6
$ws.Cells.Item(19, 6).Value = 874
$ws.Cells.Item(20, 1).Value = 17
$ws.Cells.Item(20, 3).Value = 'Greek Islands'
$ws.Cells.Item(20, 4).Value = 2
$ws.Cells.Item(20, 5).Value = 4.6
$ws.Cells.Item(20, 6).Value = 866
$ws.Cells.Item(21, 1).Value = 21
$ws.Cells.Item(21, 3).Value = 'Grindstone on the Monon'
$ws.Cells.Item(21, 4).Value = 2
$ws.Cells.Item(21, 5).Value = 4.4
$ws.Cells.Item(21, 6).Value = 554
$ws.Cells.Item(22, 1).Value = 49
$ws.Cells.Item(22, 3).Value = 'His Place Eatery - Chicken & Waffles, Ribs and Soul Food'
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(22, 5).Value = 4.5
$ws.Cells.Item(22, 6).Value = 2164
$ws.Cells.Item(23, 1).Value = 53
$ws.Cells.Item(23, 3).Value = 'IHOP'
$ws.Cells.Item(23, 4).Value = 1
$ws.Cells.Item(23, 5).Value = 4.1
$ws.Cells.Item(23, 6).Value = 2047
$ws.Cells.Item(24, 1).Value = 52
$ws.Cells.Item(24, 3).Value = 'Iaria''s Italian Restaurant'
$ws.Cells.Item(24, 4).Value = 2
$ws.Cells.Item(24, 5).Value = 4.6
$ws.Cells.Item(24, 6).Value = 1133
$ws.Cells.Item(25, 1).Value = 20
$ws.Cells.Item(25, 3).Value = 'Iron Skillet Restaurant'
$ws.Cells.Item(25, 4).Value = 2
$ws.Cells.Item(25, 5).Value = 4.5
$ws.Cells.Item(25, 6).Value = 470
$ws.Cells.Item(26, 1).Value = 14
$ws.Cells.Item(26, 3).Value = 'Livery'
$ws.Cells.Item(26, 4).Value = 2
$ws.Cells.Item(26, 5).Value = 4.7
$ws.Cells.Item(26, 6).Value = 1490
$ws.Cells.Item(27, 1).Value = 9
$ws.Cells.Item(27, 3).Value = 'Maggiano''s Little Italy'
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = 4.4
$ws.Cells.Item(27, 6).Value = 2257
$ws.Cells.Item(28, 1).Value = 6
$ws.Cells.Item(28, 3).Value = 'Major Restaurant'
$ws.Cells.Item(28, 4).Value = 2
$ws.Cells.Item(28, 5).Value = 4.6
$ws.Cells.Item(28, 6).Value = 427
$ws.Cells.Item(29, 1).Value = 45
$ws.Cells.Item(29, 3).Value = 'Mama Carolla''s'
$ws.Cells.Item(29, 4).Value = 2
$ws.Cells.Item(29, 5).Value = 4.7
$ws.Cells.Item(29, 6).Value = 1639
$ws.Cells.Item(30, 1).Value = 8
$ws.Cells.Item(30, 3).Value = 'Meridian Restaurant & Bar'
$ws.Cells.Item(30, 4).Value = 3
$ws.Cells.Item(30, 5).Value = 4.5
$ws.Cells.Item(30, 6).Value = 365
$ws.Cells.Item(31, 1).Value = 33
$ws.Cells.Item(31, 3).Value = 'Mesh'
$ws.Cells.Item(31, 4).Value = 3
$ws.Cells.Item(31, 5).Value = 4.4
$ws.Cells.Item(31, 6).Value = 1150
$ws.Cells.Item(32, 1).Value = 42
$ws.Cells.Item(32, 3).Value = 'Mimi Blue Restaurants'
$ws.Cells.Item(32, 4).Value = 2
$ws.Cells.Item(32, 5).Value = 4.5
$ws.Cells.Item(32, 6).Value = 666
$ws.Cells.Item(33, 1).Value = 38
$ws.Cells.Item(33, 3).Value = 'Nada'
$ws.Cells.Item(33, 4).Value = 2
$ws.Cells.Item(33, 5).Value = 4.4
$ws.Cells.Item(33, 6).Value = 1952
$ws.Cells.Item(34, 1).Value = 24
$ws.Cells.Item(34, 3).Value = 'Nesso'
$ws.Cells.Item(34, 4).Value = ""
$ws.Cells.Item(34, 5).Value = 4.7
$ws.Cells.Item(34, 6).Value = 219
$ws.Cells.Item(35, 1).Value = 36
$ws.Cells.Item(35, 3).Value = 'Ocean Prime'
$ws.Cells.Item(35, 4).Value = 4
$ws.Cells.Item(35, 5).Value = 4.6
$ws.Cells.Item(35, 6).Value = 958
$ws.Cells.Item(36, 1).Value = 54
$ws.Cells.Item(36, 3).Value = 'Olive Garden Italian Restaurant'
$ws.Cells.Item(36, 4).Value = 2
$ws.Cells.Item(36, 5).Value = 4.2
$ws.Cells.Item(36, 6).Value = 2759
$ws.Cells.Item(37, 1).Value = 39
$ws.Cells.Item(37, 3).Value = 'Pasto Italiano Restaurant & Bar'
$ws.Cells.Item(37, 4).Value = 2
$ws.Cells.Item(37, 5).Value = 4.7
$ws.Cells.Item(37, 6).Value = 195
$ws.Cells.Item(38, 1).Value = 15
$ws.Cells.Item(38, 3).Value = 'Ristorante Roma'
$ws.Cells.Item(38, 4).Value = ""
$ws.Cells.Item(38, 5).Value = 4.7
$ws.Cells.Item(38, 6).Value = 159
$ws.Cells.Item(39, 1).Value = 3
$ws.Cells.Item(39, 3).Value = 'Rusty Bucket Restaurant and Tavern'
$ws.Cells.Item(39, 4).Value = 2
$ws.Cells.Item(39, 5).Value = 4.4
$ws.Cells.Item(39, 6).Value = 946
$ws.Cells.Item(40, 1).Value = 19
$ws.Cells.Item(40, 3).Value = 'Ruth''s Chris Steak House'
$ws.Cells.Item(40, 4).Value = 4
$ws.Cells.Item(40, 5).Value = 4.6
$ws.Cells.Item(40, 6).Value = 969
$ws.Cells.Item(41, 1).Value = 16
$ws.Cells.Item(41, 3).Value = 'Sahm''s Restaurant'
$ws.Cells.Item(41, 4).Value = 2
$ws.Cells.Item(41, 5).Value = 4.5
$ws.Cells.Item(41, 6).Value = 793
$ws.Cells.Item(42, 1).Value = 2
$ws.Cells.Item(42, 3).Value = 'Seasons 52'
$ws.Cells.Item(42, 4).Value = 2
$ws.Cells.Item(42, 5).Value = 4.5
$ws.Cells.Item(42, 6).Value = 1339
$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 3).Value = 'Slapfish'
$ws.Cells.Item(43, 4).Value = 2
$ws.Cells.Item(43, 5).Value = 4.6
$ws.Cells.Item(43, 6).Value = 317
$ws.Cells.Item(44, 1).Value = 1
$ws.Cells.Item(44, 3).Value = 'The Capital Grille'
$ws.Cells.Item(44, 4).Value = 4
$ws.Cells.Item(44, 5).Value = 4.6
$ws.Cells.Item(44, 6).Value = 821
$ws.Cells.Item(45, 1).Value = 34
$ws.Cells.Item(45, 3).Value = 'The Cheesecake Factory'
$ws.Cells.Item(45, 4).Value = 2
$ws.Cells.Item(45, 5).Value = 4.2
$ws.Cells.Item(45, 6).Value = 3306
$ws.Cells.Item(46, 1).Value = 50
$ws.Cells.Item(46, 3).Value = 'The Italian House on Park'
$ws.Cells.Item(46, 4).Value = 2
$ws.Cells.Item(46, 5).Value = 4.8
$ws.Cells.Item(46, 6).Value = 544
$ws.Cells.Item(47, 1).Value = 57
$ws.Cells.Item(47, 3).Value = 'The Oceanaire Seafood Room'
$ws.Cells.Item(47, 4).Value = 3
$ws.Cells.Item(47, 5).Value = 4.5
$ws.Cells.Item(47, 6).Value = 975
$ws.Cells.Item(48, 1).Value = 51
$ws.Cells.Item(48, 3).Value = 'The Rathskeller'
$ws.Cells.Item(48, 4).Value = 2
$ws.Cells.Item(48, 5).Value = 4.5
$ws.Cells.Item(48, 6).Value = 2685
$ws.Cells.Item(49, 1).Value = 58
$ws.Cells.Item(49, 3).Value = 'Tijuana Flats - Fishers'
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(49, 5).Value = 4.6
$ws.Cells.Item(49, 6).Value = 1460
$ws.Cells.Item(50, 1).Value = 13
$ws.Cells.Item(50, 3).Value = 'Tinker Street Restaurant'
$ws.Cells.Item(50, 4).Value = 3
$ws.Cells.Item(50, 5).Value = 4.7
$ws.Cells.Item(50, 6).Value = 702
$ws.Cells.Item(51, 1).Value = 46
$ws.Cells.Item(51, 3).Value = 'Twin Peaks Restaurant'
$ws.Cells.Item(51, 4).Value = 2
$ws.Cells.Item(51, 5).Value = 4.5
$ws.Cells.Item(51, 6).Value = 3906
$ws.Cells.Item(52, 1).Value = 32
$ws.Cells.Item(52, 3).Value = 'Weber Grill Restaurant'
$ws.Cells.Item(52, 4).Value = 2
$ws.Cells.Item(52, 5).Value = 4.2
$ws.Cells.Item(52, 6).Value = 2307
$ws.Cells.Item(53, 1).Value = 7
$ws.Cells.Item(53, 3).Value = 'Yard House'
$ws.Cells.Item(53, 4).Value = 2
$ws.Cells.Item(53, 5).Value = 4.4
$ws.Cells.Item(53, 6).Value = 2326
$ws.Cells.Item(54, 1).Value = 11
$ws.Cells.Item(54, 3).Value = 'Yats'
$ws.Cells.Item(54, 4).Value = 1
$ws.Cells.Item(54, 5).Value = 4.8
$ws.Cells.Item(54, 6).Value = 1279
